# Generate Report for handback
#
# For each locale sheet (zh-cn, de-de):
#   - Status column (B) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for the two real rows (2 and 3).
#   - Two new columns get populated now that handback has happened:
#       E = Latest Target File    (same file/link as row 2's column A)
#       F = Latest Handback File  (same file/link as row 2's column C)
#   - Column G (Latest Handback DateTime) gets the real handback
#     timestamp instead of the 0001-01-01 00:00:00 placeholder.
#   - The hyperlink style (underline + the workbook's custom blue) used
#     by the existing A/C hyperlink cells is re-applied to the new E/F
#     cells.

$hyperlinkColor = 15570276   # matches the workbook's custom "HyperLink" style (RGB 0x6495ED)

function Apply-HandbackRow {
    param($ws, $row, $srcRow, $handbackDateTime)

    # Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
    $ws.Cells.Item($row, 2).Value = "Handed back: in sync with en-US"

    # Source file (column A) + its hyperlink target, read from $srcRow
    # (row 3's target/handback files mirror row 2's, since the second
    # file depends on / reuses the first file's translation).
    $fileName = $ws.Cells.Item($srcRow, 1).Text
    $fileUrl = $null
    $xlfName = $ws.Cells.Item($srcRow, 3).Text
    $xlfUrl = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq $srcRow -and $hl.Range.Column -eq 1) {
            $fileUrl = $hl.Address
        }
        if ($hl.Range.Row -eq $srcRow -and $hl.Range.Column -eq 3) {
            $xlfUrl = $hl.Address
        }
    }

    # E = Latest Target File (same file/link as the Source File Name column)
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.Value = $fileName
    $ws.Hyperlinks.Add($eCell, $fileUrl, "", "", $fileName) | Out-Null
    $eCell.Font.Underline = 2
    $eCell.Font.Color = $hyperlinkColor

    # F = Latest Handback File (same file/link as the Latest Handoff File column)
    $fCell = $ws.Cells.Item($row, 6)
    $fCell.Value = $xlfName
    $ws.Hyperlinks.Add($fCell, $xlfUrl, "", "", $xlfName) | Out-Null
    $fCell.Font.Underline = 2
    $fCell.Font.Color = $hyperlinkColor

    # G = Latest Handback DateTime (real timestamp instead of the 0001-01-01 placeholder)
    $ws.Cells.Item($row, 7).Value = $handbackDateTime
}

$wb = $excel.ActiveWorkbook

# The Overview sheet mirrors each locale's status in its own columns via
# the same shared string, so it flips to "Handed back..." as well.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 2).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(2, 3).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3, 2).Value = "Handed back: in sync with en-US"
$wsOverview.Cells.Item(3, 3).Value = "Handed back: in sync with en-US"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Apply-HandbackRow $wsZhCn 2 2 "2016-01-21 02:23:27"
Apply-HandbackRow $wsZhCn 3 2 "2016-01-21 02:23:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
Apply-HandbackRow $wsDeDe 2 2 "2016-01-21 02:23:46"
Apply-HandbackRow $wsDeDe 3 2 "2016-01-21 02:23:46"
